$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.234.70"
$ws.Range("E2").Value = "  +0.58%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.842.13"
$ws.Range("E3").Value = "  +0.33%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.83%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.90"
$ws.Range("E5").Value = "  +1.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6316"
$ws.Range("E6").Value = "  +0.48%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.008"
$ws.Range("E7").Value = "  +0.73%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07502"
$ws.Range("E8").Value = "  -2.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2940"
$ws.Range("E9").Value = "  +0.59%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.20"
$ws.Range("E10").Value = "  +1.81%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07765"
$ws.Range("E11").Value = "  +0.25%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.839.10"
$ws.Range("E12").Value = "  +0.21%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.006"
$ws.Range("E13").Value = "  +0.88%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6711"
$ws.Range("E14").Value = "  +0.67%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "83.34"
$ws.Range("E15").Value = "  +0.49%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009378"
$ws.Range("E16").Value = "  -9.08%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.074"
$ws.Range("E17").Value = "  +0.97%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.239.50"
$ws.Range("E18").Value = "  +0.57%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.66"
$ws.Range("E19").Value = "  +2.41%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "224.72"
$ws.Range("E20").Value = "  -0.51%  "

$ws.Range("E21").Value = "  +0.90%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.162"
$ws.Range("E22").Value = "  -0.83%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "161.20"
$ws.Range("E24").Value = "  +1.74%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1407"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.545"
$ws.Range("E26").Value = "  +1.06%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.00"
$ws.Range("E27").Value = "  +0.36%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.509"
$ws.Range("E28").Value = "  +1.14%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.171"
$ws.Range("E29").Value = "  +2.46%  "

$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05555"
$ws.Range("E30").Value = "  +6.61%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.088"
$ws.Range("E31").Value = "  +1.31%  "

$ws.Range("E32").Value = "  +0.77%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7543"
$ws.Range("E33").Value = "  +2.10%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.863"
$ws.Range("E34").Value = "  +0.58%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.141"
$ws.Range("E35").Value = "  -0.40%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.624"
$ws.Range("E36").Value = "  -2.80%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.239.51"
$ws.Range("E37").Value = "  -2.35%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.764"
$ws.Range("E38").Value = "  +0.05%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01792"
$ws.Range("E39").Value = "  +0.38%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.601"
$ws.Range("E40").Value = "  +4.11%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8961"
$ws.Range("E41").Value = "  -0.16%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.008"
$ws.Range("E42").Value = "  +0.65%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "102.16"
$ws.Range("E43").Value = "  +0.48%  "

$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.983.12"
$ws.Range("E44").Value = "  +0.17%  "

$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000125"
$ws.Range("E45").Value = "  +1.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.94"
$ws.Range("E46").Value = "  +2.60%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.07760"
$ws.Range("E47").Value = "  +13.12%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5119"
$ws.Range("E48").Value = "  +0.04%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4075"
$ws.Range("E49").Value = "  +1.50%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.079"
$ws.Range("E50").Value = "  +2.65%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05822"
$ws.Range("E51").Value = "  +1.01%  "
